$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 13.45843057308089
$ws.Range("D2").Value = 5.67360955427247
$ws.Range("E2").Value = 13.82782050721507
$ws.Range("F2").Value = 34.00327528687026
$ws.Range("G2").Value = 3.640837487773293
$ws.Range("L2").Value = 9.225315554510207
$ws.Range("N2").Value = 18.99769123883691
$ws.Range("O2").Value = 29.72179095132885
$ws.Range("C3").Value = 13.3661654808872
$ws.Range("D3").Value = 5.686663600237641
$ws.Range("E3").Value = 13.77440445474033
$ws.Range("F3").Value = 33.4555667402718
$ws.Range("G3").Value = 3.645260097764718
$ws.Range("L3").Value = 9.233124427980423
$ws.Range("N3").Value = 18.40031508502705
$ws.Range("O3").Value = 29.35654874916323
$ws.Range("C4").Value = 13.31305664959141
$ws.Range("D4").Value = 5.695707863816937
$ws.Range("E4").Value = 13.74516280287087
$ws.Range("F4").Value = 33.1261047885181
$ws.Range("G4").Value = 3.648115624878236
$ws.Range("L4").Value = 9.239670170994417
$ws.Range("N4").Value = 18.02485520896364
$ws.Range("O4").Value = 29.13942427315438
$ws.Range("C5").Value = 13.29232169284367
$ws.Range("D5").Value = 5.699651528292453
$ws.Range("E5").Value = 13.73414686557506
$ws.Range("F5").Value = 32.99374231114239
$ws.Range("G5").Value = 3.649314633341758
$ws.Range("L5").Value = 9.242777148845201
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 29.05283770159657
$ws.Range("C6").Value = 13.28893392007817
$ws.Range("D6").Value = 5.700321935737549
$ws.Range("E6").Value = 13.73237220620931
$ws.Range("F6").Value = 32.97188325160016
$ws.Range("G6").Value = 3.649515867633657
$ws.Range("L6").Value = 9.243319581575312
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 29.03857723123535
$ws.Range("C7").Value = 13.31277331633427
$ws.Range("D7").Value = 5.695760005535817
$ws.Range("E7").Value = 13.74501058573987
$ws.Range("F7").Value = 33.12431179315897
$ws.Range("G7").Value = 3.648131651765919
$ws.Range("L7").Value = 9.239710294240801
$ws.Range("N7").Value = 18.02277304766463
$ws.Range("O7").Value = 29.13824874237422
$ws.Range("C8").Value = 13.42589282295319
$ws.Range("D8").Value = 5.677896594854621
$ws.Range("E8").Value = 13.80866860424339
$ws.Range("F8").Value = 33.81311469682657
$ws.Range("G8").Value = 3.642333429923539
$ws.Range("L8").Value = 9.227644183275331
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("O8").Value = 29.59443518759517
$ws.Range("C9").Value = 13.67501977545935
$ws.Range("D9").Value = 5.651062814987768
$ws.Range("E9").Value = 13.96139696018238
$ws.Range("F9").Value = 35.21014331858598
$ws.Range("G9").Value = 3.632067476928689
$ws.Range("L9").Value = 9.217907576451738
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 30.54108561838983
$ws.Range("C10").Value = 13.87351330624774
$ws.Range("D10").Value = 5.636387773005403
$ws.Range("E10").Value = 14.09010922011097
$ws.Range("F10").Value = 36.25410360309793
$ws.Range("G10").Value = 3.625188936927255
$ws.Range("L10").Value = 9.219279807822959
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 31.26204015738269
$ws.Range("C11").Value = 13.96689110188261
$ws.Range("D11").Value = 5.630815051330783
$ws.Range("E11").Value = 14.15211736002805
$ws.Range("F11").Value = 36.73061710149484
$ws.Range("G11").Value = 3.622201851554292
$ws.Range("L11").Value = 9.221760514494774
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 31.59417816948042
$ws.Range("C12").Value = 14.00266896467751
$ws.Range("D12").Value = 5.628864128346152
$ws.Range("E12").Value = 14.17608267467964
$ws.Range("F12").Value = 36.91111181689754
$ws.Range("G12").Value = 3.621090983327357
$ws.Range("L12").Value = 9.222966993544306
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 31.72043532706188
$ws.Range("C13").Value = 13.99494539814866
$ws.Range("D13").Value = 5.629277196090541
$ws.Range("E13").Value = 14.17089999553264
$ws.Range("F13").Value = 36.8722396914956
$ws.Range("G13").Value = 3.621329329158419
$ws.Range("L13").Value = 9.222695277311612
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 31.69322391484196
$ws.Range("C14").Value = 13.96982634539408
$ws.Range("D14").Value = 5.630651350069042
$ws.Range("E14").Value = 14.15407937401407
$ws.Range("F14").Value = 36.74546654308784
$ws.Range("G14").Value = 3.622110054209239
$ws.Range("L14").Value = 9.221854419159639
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 31.60455638953862
$ws.Range("C15").Value = 13.95449382897115
$ws.Range("D15").Value = 5.631513831998677
$ws.Range("E15").Value = 14.14383891528427
$ws.Range("F15").Value = 36.66781545417744
$ws.Range("G15").Value = 3.622590907166923
$ws.Range("L15").Value = 9.221374153538857
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 31.55030453097894
$ws.Range("C16").Value = 13.86747070858826
$ws.Range("D16").Value = 5.636774222147579
$ws.Range("E16").Value = 14.08612531612221
$ws.Range("F16").Value = 36.2229790292487
$ws.Range("G16").Value = 3.625386995109486
$ws.Range("L16").Value = 9.219155059317094
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 31.24040815154042
$ws.Range("C17").Value = 13.81485638976041
$ws.Range("D17").Value = 5.640284338997454
$ws.Range("E17").Value = 14.05159630043378
$ws.Range("F17").Value = 35.95036740577823
$ws.Range("G17").Value = 3.627138571631213
$ws.Range("L17").Value = 9.21826934364171
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 31.05128235415316
$ws.Range("C18").Value = 13.78488571247737
$ws.Range("D18").Value = 5.642407044972952
$ws.Range("E18").Value = 14.03206227613792
$ws.Range("F18").Value = 35.79373059941134
$ws.Range("G18").Value = 3.628159406072558
$ws.Range("L18").Value = 9.217934644857733
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 30.94290182568167
$ws.Range("C19").Value = 13.77478901584949
$ws.Range("D19").Value = 5.643143558397747
$ws.Range("E19").Value = 14.02550478367911
$ws.Range("F19").Value = 35.7407295693666
$ws.Range("G19").Value = 3.628507344416974
$ws.Range("L19").Value = 9.217851327196907
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 30.90627823508273
$ws.Range("C20").Value = 13.82042726717809
$ws.Range("D20").Value = 5.639899934923041
$ws.Range("E20").Value = 14.05523830757427
$ws.Range("F20").Value = 35.97937193093676
$ws.Range("G20").Value = 3.626950730202652
$ws.Range("L20").Value = 9.218345543099945
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("O20").Value = 31.0713746036034
$ws.Range("C21").Value = 13.977193290522
$ws.Range("D21").Value = 5.630243397387904
$ws.Range("E21").Value = 14.15900696467714
$ws.Range("F21").Value = 36.78270297203284
$ws.Range("G21").Value = 3.621880187184309
$ws.Range("L21").Value = 9.222094150543212
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 31.6305879843048
$ws.Range("C22").Value = 14.0820713647514
$ws.Range("D22").Value = 5.624861408286409
$ws.Range("E22").Value = 14.22964135465732
$ws.Range("F22").Value = 37.30793404426885
$ws.Range("G22").Value = 3.618684413146165
$ws.Range("L22").Value = 9.226100839348996
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 31.99883269342842
$ws.Range("C23").Value = 14.02588293528818
$ws.Range("D23").Value = 5.627648625605701
$ws.Range("E23").Value = 14.19168929947749
$ws.Range("F23").Value = 37.02764794139411
$ws.Range("G23").Value = 3.620379297130368
$ws.Range("L23").Value = 9.223819945605371
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 31.80207850291564
$ws.Range("C24").Value = 13.81790780568118
$ws.Range("D24").Value = 5.640073397964894
$ws.Range("E24").Value = 14.05359076728055
$ws.Range("F24").Value = 35.96625868700839
$ws.Range("G24").Value = 3.627035610225261
$ws.Range("L24").Value = 9.218310549773868
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 31.06228979734288
$ws.Range("C25").Value = 13.60481832645839
$ws.Range("D25").Value = 5.657440090125234
$ws.Range("E25").Value = 13.91714100409435
$ws.Range("F25").Value = 34.82835901985329
$ws.Range("G25").Value = 3.634727432625657
$ws.Range("L25").Value = 9.219046490638688
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 30.28005906737457
